# Review-Log-Sheet.xlsx — "adding rev-12 to review log related to the
# integration testing"
#
# 1. Row 12 ("rev-10") Status cell (G12) flips from PENDING -> DONE.
# 2. A brand-new review row (row 14, "rev-12" / Integration Test) is
#    appended, re-using the same formatting as row 12 (border/fill/
#    alignment) and the row-12 PENDING status look.
# 3. The sheet's active selection moves from F13 to G13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 2. Build new row 14 first, while row 12 still carries the PENDING
#        (G-column) look we want row 14's Status cell to copy. ---------
$null = $ws.Range("A12:H12").Copy()
$null = $ws.Range("A14:H14").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(14).RowHeight = 45

$ws.Range("A14").Value = "rev-12"
$ws.Range("B14").Value = "Integration Test "
$ws.Range("C14").Value = "bad test scenarios. More scenarios needed with each two modules"
$ws.Range("D14").Value = 42586
$ws.Range("E14").Value = "Med"
$ws.Range("F14").Value = "Walaa"
$ws.Range("G14").Value = "PENDING"
$ws.Range("H14").Value = "seif eldin"

# --- 1. Now flip row 12's Status (G12) to DONE, copying G11's DONE
#        formatting (green fill) onto it before overwriting the text. ---
$null = $ws.Range("G11").Copy()
$null = $ws.Range("G12").PasteSpecial(-4122)       # xlPasteFormats
$ws.Range("G12").Value = "DONE"

# --- 3. Move the active selection to G13. ------------------------------
$null = $ws.Range("G13").Select()
